$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PrimaryContact column (D) is becoming a boolean flag: every contact row
# below the first (D2, already 1) is normalized to 1.
$ws.Range("D3:D11").Value = 1

# Update the active cell / selection to match the saved view state.
$ws.Range("D14").Select()
